# Update Maria's "Idade" (age) from 13 to 14.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B3").Value = 14
